$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 9).Value = 104
$ws.Cells.Item(3, 7).Value = 115
$ws.Cells.Item(3, 9).Value = 69
$ws.Cells.Item(4, 7).Value = 14
$ws.Cells.Item(5, 7).Value = 14
$ws.Cells.Item(6, 7).Value = 15
$ws.Cells.Item(6, 9).Value = 77
$ws.Cells.Item(7, 7).Value = 74
$ws.Cells.Item(7, 9).Value = 18
$ws.Cells.Item(8, 7).Value = 80
$ws.Cells.Item(8, 9).Value = 12
$ws.Cells.Item(9, 7).Value = 3
$ws.Cells.Item(9, 9).Value = 89
$ws.Cells.Item(10, 7).Value = 174
$ws.Cells.Item(10, 9).Value = 10
$ws.Cells.Item(11, 9).Value = 50
$ws.Cells.Item(12, 7).Value = 41
$ws.Cells.Item(12, 9).Value = 141
$ws.Cells.Item(13, 7).Value = 17
$ws.Cells.Item(13, 9).Value = 166
$ws.Cells.Item(14, 7).Value = 17
$ws.Cells.Item(14, 9).Value = 166
$ws.Cells.Item(15, 7).Value = 130
$ws.Cells.Item(15, 9).Value = 53
$ws.Cells.Item(16, 7).Value = 50
$ws.Cells.Item(16, 9).Value = 42
$ws.Cells.Item(17, 7).Value = 50
$ws.Cells.Item(17, 9).Value = 41
$ws.Cells.Item(18, 7).Value = 12
$ws.Cells.Item(18, 9).Value = 80
$ws.Cells.Item(19, 7).Value = 7
$ws.Cells.Item(19, 9).Value = 85
$ws.Cells.Item(20, 7).Value = 53
$ws.Cells.Item(20, 9).Value = 38
$ws.Cells.Item(21, 7).Value = 49
$ws.Cells.Item(21, 9).Value = 42
$ws.Cells.Item(22, 7).Value = 49
$ws.Cells.Item(22, 9).Value = 42
$ws.Cells.Item(23, 9).Value = 9
$ws.Cells.Item(24, 9).Value = 9
$ws.Cells.Item(25, 7).Value = 32
$ws.Cells.Item(26, 7).Value = 76
$ws.Cells.Item(26, 9).Value = 108
$ws.Cells.Item(27, 7).Value = 118
$ws.Cells.Item(27, 9).Value = 66
$ws.Cells.Item(28, 7).Value = 151
$ws.Cells.Item(28, 9).Value = 32
$ws.Cells.Item(29, 7).Value = 53
$ws.Cells.Item(29, 9).Value = 38
$ws.Cells.Item(30, 7).Value = 12
$ws.Cells.Item(30, 9).Value = 80
$ws.Cells.Item(31, 7).Value = 83
$ws.Cells.Item(31, 9).Value = 9
$ws.Cells.Item(32, 7).Value = 81
$ws.Cells.Item(32, 9).Value = 11
$ws.Cells.Item(33, 9).Value = 19
$ws.Cells.Item(34, 7).Value = 60
$ws.Cells.Item(34, 9).Value = 31
$ws.Cells.Item(35, 7).Value = 14
$ws.Cells.Item(35, 9).Value = 78
$ws.Cells.Item(36, 7).Value = 29
$ws.Cells.Item(36, 9).Value = 63
$ws.Cells.Item(37, 7).Value = 47
$ws.Cells.Item(37, 9).Value = 44
$ws.Cells.Item(38, 7).Value = 129
$ws.Cells.Item(38, 9).Value = 54
$ws.Cells.Item(39, 7).Value = 173
$ws.Cells.Item(39, 9).Value = 11
$ws.Cells.Item(40, 7).Value = 59
$ws.Cells.Item(40, 9).Value = 123
$ws.Cells.Item(41, 7).Value = 129
$ws.Cells.Item(41, 9).Value = 54
$ws.Cells.Item(42, 9).Value = 29
$ws.Cells.Item(43, 7).Value = 138
$ws.Cells.Item(43, 9).Value = 45
$ws.Cells.Item(44, 7).Value = 137
$ws.Cells.Item(44, 9).Value = 46
$ws.Cells.Item(45, 7).Value = 130
$ws.Cells.Item(45, 9).Value = 53
$ws.Cells.Item(46, 7).Value = 130
$ws.Cells.Item(46, 9).Value = 53
$ws.Cells.Item(47, 7).Value = 116
$ws.Cells.Item(47, 9).Value = 68
$ws.Cells.Item(48, 7).Value = 118
$ws.Cells.Item(48, 9).Value = 66
$ws.Cells.Item(49, 7).Value = 84
$ws.Cells.Item(49, 9).Value = 100
$ws.Cells.Item(50, 7).Value = 11
$ws.Cells.Item(50, 9).Value = 172
$ws.Cells.Item(51, 7).Value = 114
$ws.Cells.Item(51, 9).Value = 70
$ws.Cells.Item(52, 7).Value = 97
$ws.Cells.Item(52, 9).Value = 87
$ws.Cells.Item(53, 7).Value = 137
$ws.Cells.Item(53, 9).Value = 46
$ws.Cells.Item(54, 7).Value = 40
$ws.Cells.Item(54, 9).Value = 51
$ws.Cells.Item(55, 7).Value = 9
$ws.Cells.Item(56, 7).Value = 85
$ws.Cells.Item(56, 9).Value = 7
$ws.Cells.Item(57, 7).Value = 40
$ws.Cells.Item(57, 9).Value = 51
$ws.Cells.Item(58, 7).Value = 167
$ws.Cells.Item(58, 9).Value = 17
$ws.Cells.Item(59, 7).Value = 45
$ws.Cells.Item(59, 9).Value = 46
$ws.Cells.Item(60, 7).Value = 31
$ws.Cells.Item(61, 7).Value = 17
$ws.Cells.Item(61, 9).Value = 77
$ws.Cells.Item(62, 7).Value = 178
$ws.Cells.Item(62, 9).Value = 6
$ws.Cells.Item(63, 7).Value = 42
$ws.Cells.Item(63, 9).Value = 140
$ws.Cells.Item(64, 7).Value = 28
$ws.Cells.Item(65, 7).Value = 28
$ws.Cells.Item(66, 7).Value = 147
$ws.Cells.Item(66, 9).Value = 36
$ws.Cells.Item(67, 7).Value = 58
$ws.Cells.Item(67, 9).Value = 33
$ws.Cells.Item(68, 7).Value = 40
$ws.Cells.Item(68, 9).Value = 51
$ws.Cells.Item(69, 7).Value = 87
$ws.Cells.Item(69, 9).Value = 5
$ws.Cells.Item(70, 7).Value = 31
$ws.Cells.Item(70, 9).Value = 60
$ws.Cells.Item(71, 7).Value = 31
$ws.Cells.Item(72, 7).Value = 29
$ws.Cells.Item(72, 9).Value = 63
$ws.Cells.Item(73, 7).Value = 47
$ws.Cells.Item(73, 9).Value = 44
$ws.Cells.Item(74, 7).Value = 15
$ws.Cells.Item(74, 9).Value = 168
$ws.Cells.Item(75, 7).Value = 119
$ws.Cells.Item(75, 9).Value = 65
$ws.Cells.Item(76, 7).Value = 76
$ws.Cells.Item(76, 9).Value = 108
$ws.Cells.Item(77, 7).Value = 50
$ws.Cells.Item(77, 9).Value = 132
$ws.Cells.Item(78, 7).Value = 10
$ws.Cells.Item(78, 9).Value = 173
$ws.Cells.Item(79, 7).Value = 157
$ws.Cells.Item(79, 9).Value = 26
$ws.Cells.Item(80, 9).Value = 18
$ws.Cells.Item(81, 7).Value = 132
$ws.Cells.Item(81, 9).Value = 51
$ws.Cells.Item(82, 9).Value = 52
$ws.Cells.Item(83, 7).Value = 61
$ws.Cells.Item(83, 9).Value = 121
$ws.Cells.Item(84, 7).Value = 182
$ws.Cells.Item(84, 9).Value = 2
$ws.Cells.Item(85, 7).Value = 180
$ws.Cells.Item(85, 9).Value = 4
$ws.Cells.Item(86, 7).Value = 151
$ws.Cells.Item(86, 9).Value = 32
$ws.Cells.Item(87, 7).Value = 42
$ws.Cells.Item(87, 9).Value = 324
$ws.Cells.Item(88, 7).Value = 133
$ws.Cells.Item(88, 9).Value = 50
$ws.Cells.Item(89, 7).Value = 133
$ws.Cells.Item(89, 9).Value = 50
$ws.Cells.Item(90, 7).Value = 44
$ws.Cells.Item(90, 9).Value = 138
$ws.Cells.Item(91, 7).Value = 136
$ws.Cells.Item(91, 9).Value = 51
$ws.Cells.Item(92, 7).Value = 49
$ws.Cells.Item(92, 9).Value = 133
$ws.Cells.Item(93, 7).Value = 132
$ws.Cells.Item(93, 9).Value = 51
$ws.Cells.Item(94, 7).Value = 132
$ws.Cells.Item(94, 9).Value = 51
$ws.Cells.Item(95, 7).Value = 3
$ws.Cells.Item(95, 9).Value = 179
$ws.Cells.Item(96, 7).Value = 114
$ws.Cells.Item(96, 9).Value = 70
$ws.Cells.Item(97, 7).Value = 95
$ws.Cells.Item(97, 9).Value = 89
$ws.Cells.Item(98, 7).Value = 42
$ws.Cells.Item(98, 9).Value = 140
$ws.Cells.Item(99, 7).Value = 143
$ws.Cells.Item(99, 9).Value = 40
$ws.Cells.Item(100, 7).Value = 113
$ws.Cells.Item(100, 9).Value = 71
$ws.Cells.Item(101, 7).Value = 38
$ws.Cells.Item(101, 9).Value = 144
$ws.Cells.Item(102, 7).Value = 143
$ws.Cells.Item(102, 9).Value = 40
$ws.Cells.Item(103, 7).Value = 114
$ws.Cells.Item(103, 9).Value = 70
$ws.Cells.Item(104, 7).Value = 114
$ws.Cells.Item(104, 9).Value = 70
$ws.Cells.Item(105, 7).Value = 2
$ws.Cells.Item(106, 7).Value = 151
$ws.Cells.Item(106, 9).Value = 32
$ws.Cells.Item(107, 7).Value = 51
$ws.Cells.Item(107, 9).Value = 40
$ws.Cells.Item(108, 7).Value = 51
$ws.Cells.Item(108, 9).Value = 40
$ws.Cells.Item(109, 7).Value = 51
$ws.Cells.Item(109, 9).Value = 40
$ws.Cells.Item(110, 7).Value = 88
$ws.Cells.Item(110, 9).Value = 4
$ws.Cells.Item(111, 7).Value = 25
$ws.Cells.Item(111, 9).Value = 67
$ws.Cells.Item(112, 7).Value = 64
$ws.Cells.Item(112, 9).Value = 27
$ws.Cells.Item(113, 7).Value = 30
$ws.Cells.Item(113, 9).Value = 153
$ws.Cells.Item(114, 7).Value = 33
$ws.Cells.Item(114, 9).Value = 59
$ws.Cells.Item(115, 7).Value = 18
$ws.Cells.Item(115, 9).Value = 74
$ws.Cells.Item(116, 7).Value = 75
$ws.Cells.Item(116, 9).Value = 17
$ws.Cells.Item(117, 7).Value = 121
$ws.Cells.Item(117, 9).Value = 63
$ws.Cells.Item(118, 7).Value = 142
$ws.Cells.Item(118, 9).Value = 41
$ws.Cells.Item(119, 7).Value = 132
$ws.Cells.Item(119, 9).Value = 51
$ws.Cells.Item(120, 7).Value = 132
$ws.Cells.Item(120, 9).Value = 51
$ws.Cells.Item(121, 7).Value = 132
$ws.Cells.Item(121, 9).Value = 51
$ws.Cells.Item(122, 7).Value = 132
$ws.Cells.Item(122, 9).Value = 51
$ws.Cells.Item(123, 7).Value = 37
$ws.Cells.Item(123, 9).Value = 55
$ws.Cells.Item(124, 7).Value = 39
$ws.Cells.Item(124, 9).Value = 143
$ws.Cells.Item(125, 7).Value = 23
$ws.Cells.Item(125, 9).Value = 69
$ws.Cells.Item(126, 7).Value = 10
$ws.Cells.Item(126, 9).Value = 82
$ws.Cells.Item(127, 7).Value = 24
$ws.Cells.Item(127, 9).Value = 68
$ws.Cells.Item(128, 7).Value = 18
$ws.Cells.Item(128, 9).Value = 74
$ws.Cells.Item(129, 7).Value = 78
$ws.Cells.Item(129, 9).Value = 14
$ws.Cells.Item(130, 7).Value = 65
$ws.Cells.Item(130, 9).Value = 26
$ws.Cells.Item(131, 7).Value = 33
$ws.Cells.Item(131, 9).Value = 59
$ws.Cells.Item(132, 7).Value = 92
$ws.Cells.Item(132, 8).Value = 45328
$ws.Cells.Item(132, 9).Value = 92
$ws.Cells.Item(133, 7).Value = 80
$ws.Cells.Item(133, 9).Value = 12
$ws.Cells.Item(134, 7).Value = 90
$ws.Cells.Item(134, 9).Value = 2
$ws.Cells.Item(135, 7).Value = 4
$ws.Cells.Item(135, 9).Value = 88
$ws.Cells.Item(136, 7).Value = 40
$ws.Cells.Item(136, 9).Value = 51
$ws.Cells.Item(137, 7).Value = 39
$ws.Cells.Item(137, 9).Value = 52
$ws.Cells.Item(138, 7).Value = 57
$ws.Cells.Item(138, 9).Value = 34
$ws.Cells.Item(139, 7).Value = 39
$ws.Cells.Item(139, 9).Value = 52
$ws.Cells.Item(140, 7).Value = 54
$ws.Cells.Item(140, 9).Value = 37
$ws.Cells.Item(141, 7).Value = 28
$ws.Cells.Item(141, 9).Value = 64
$ws.Cells.Item(142, 7).Value = 86
$ws.Cells.Item(142, 9).Value = 6
$ws.Cells.Item(143, 9).Value = 2
$ws.Cells.Item(144, 7).Value = 65
$ws.Cells.Item(144, 9).Value = 26
$ws.Cells.Item(145, 7).Value = 42
$ws.Cells.Item(145, 9).Value = 49
$ws.Cells.Item(146, 7).Value = 25
$ws.Cells.Item(147, 7).Value = 71
$ws.Cells.Item(147, 9).Value = 21
$ws.Cells.Item(148, 7).Value = 70
$ws.Cells.Item(148, 9).Value = 22
$ws.Cells.Item(149, 7).Value = 51
$ws.Cells.Item(149, 9).Value = 40
$ws.Cells.Item(150, 7).Value = 18
$ws.Cells.Item(150, 9).Value = 165
$ws.Cells.Item(151, 7).Value = 62
$ws.Cells.Item(151, 9).Value = 120
$ws.Cells.Item(152, 7).Value = 147
$ws.Cells.Item(152, 9).Value = 36
$ws.Cells.Item(153, 7).Value = 112
$ws.Cells.Item(153, 9).Value = 72
$ws.Cells.Item(154, 7).Value = 112
$ws.Cells.Item(154, 9).Value = 72
$ws.Cells.Item(155, 7).Value = 27
$ws.Cells.Item(155, 9).Value = 156
$ws.Cells.Item(156, 7).Value = 27
$ws.Cells.Item(156, 9).Value = 156
$ws.Cells.Item(157, 7).Value = 27
$ws.Cells.Item(157, 9).Value = 156
$ws.Cells.Item(158, 7).Value = 27
$ws.Cells.Item(158, 9).Value = 156
$ws.Cells.Item(159, 7).Value = 80
$ws.Cells.Item(159, 9).Value = 11
$ws.Cells.Item(160, 7).Value = 83
$ws.Cells.Item(160, 9).Value = 9
$ws.Cells.Item(161, 7).Value = 101
$ws.Cells.Item(161, 9).Value = 83
$ws.Cells.Item(162, 7).Value = 145
$ws.Cells.Item(162, 9).Value = 38
$ws.Cells.Item(163, 7).Value = 145
$ws.Cells.Item(163, 9).Value = 38
$ws.Cells.Item(164, 7).Value = 98
$ws.Cells.Item(164, 9).Value = 86
$ws.Cells.Item(165, 7).Value = 98
$ws.Cells.Item(165, 9).Value = 86
$ws.Cells.Item(166, 7).Value = 151
$ws.Cells.Item(166, 9).Value = 32
$ws.Cells.Item(167, 7).Value = 151
$ws.Cells.Item(167, 9).Value = 32
$ws.Cells.Item(168, 9).Value = 90
$ws.Cells.Item(169, 7).Value = 63
$ws.Cells.Item(169, 9).Value = 119
$ws.Cells.Item(170, 7).Value = 180
$ws.Cells.Item(170, 9).Value = 4
$ws.Cells.Item(171, 9).Value = 2
$ws.Cells.Item(172, 7).Value = 95
$ws.Cells.Item(172, 9).Value = 89
$ws.Cells.Item(173, 7).Value = 16
$ws.Cells.Item(173, 9).Value = 167
$ws.Cells.Item(174, 7).Value = 65
$ws.Cells.Item(174, 9).Value = 26
$ws.Cells.Item(175, 7).Value = 39
$ws.Cells.Item(175, 9).Value = 52
$ws.Cells.Item(176, 7).Value = 74
$ws.Cells.Item(176, 9).Value = 18
$ws.Cells.Item(177, 7).Value = 175
$ws.Cells.Item(177, 9).Value = 9
$ws.Cells.Item(178, 7).Value = 24
$ws.Cells.Item(178, 9).Value = 68
$ws.Cells.Item(179, 7).Value = 45
$ws.Cells.Item(179, 9).Value = 46
$ws.Cells.Item(180, 7).Value = 148
$ws.Cells.Item(180, 9).Value = 35
$ws.Cells.Item(181, 7).Value = 44
$ws.Cells.Item(181, 9).Value = 47
$ws.Cells.Item(182, 7).Value = 13
$ws.Cells.Item(182, 9).Value = 79
$ws.Cells.Item(183, 7).Value = 8
$ws.Cells.Item(183, 9).Value = 84
$ws.Cells.Item(184, 7).Value = 86
$ws.Cells.Item(184, 9).Value = 6
$ws.Cells.Item(185, 7).Value = 26
$ws.Cells.Item(186, 7).Value = 15
$ws.Cells.Item(186, 9).Value = 77
$ws.Cells.Item(187, 7).Value = 47
$ws.Cells.Item(187, 9).Value = 44
$ws.Cells.Item(188, 7).Value = 43
$ws.Cells.Item(188, 9).Value = 48
$ws.Cells.Item(189, 7).Value = 71
$ws.Cells.Item(189, 9).Value = 21
$ws.Cells.Item(190, 7).Value = 34
$ws.Cells.Item(190, 9).Value = 149
$ws.Cells.Item(191, 7).Value = 155
$ws.Cells.Item(191, 9).Value = 28
$ws.Cells.Item(192, 7).Value = 155
$ws.Cells.Item(192, 9).Value = 28
$ws.Cells.Item(193, 9).Value = 83
$ws.Cells.Item(194, 9).Value = 83
$ws.Cells.Item(195, 7).Value = 6
$ws.Cells.Item(195, 9).Value = 176
$ws.Cells.Item(196, 7).Value = 6
$ws.Cells.Item(196, 9).Value = 176
$ws.Cells.Item(197, 7).Value = 22
$ws.Cells.Item(197, 9).Value = 161
$ws.Cells.Item(198, 7).Value = 22
$ws.Cells.Item(198, 9).Value = 161
$ws.Cells.Item(199, 7).Value = 22
$ws.Cells.Item(199, 9).Value = 161
$ws.Cells.Item(200, 7).Value = 22
$ws.Cells.Item(200, 9).Value = 161
$ws.Cells.Item(201, 7).Value = 169
$ws.Cells.Item(201, 9).Value = 15
$ws.Cells.Item(202, 7).Value = 169
$ws.Cells.Item(202, 9).Value = 15
$ws.Cells.Item(203, 7).Value = 169
$ws.Cells.Item(203, 9).Value = 15
$ws.Cells.Item(204, 7).Value = 169
$ws.Cells.Item(204, 9).Value = 15
$ws.Cells.Item(205, 7).Value = 58
$ws.Cells.Item(205, 9).Value = 33
$ws.Cells.Item(206, 9).Value = 16
$ws.Cells.Item(207, 7).Value = 149
$ws.Cells.Item(207, 9).Value = 34
$ws.Cells.Item(208, 7).Value = 137
$ws.Cells.Item(208, 9).Value = 46
$ws.Cells.Item(209, 7).Value = 61
$ws.Cells.Item(209, 9).Value = 120
$ws.Cells.Item(210, 7).Value = 61
$ws.Cells.Item(210, 9).Value = 120
$ws.Cells.Item(211, 7).Value = 70
$ws.Cells.Item(211, 9).Value = 114
$ws.Cells.Item(212, 7).Value = 70
$ws.Cells.Item(212, 9).Value = 114
$ws.Cells.Item(213, 7).Value = 146
$ws.Cells.Item(213, 9).Value = 37
$ws.Cells.Item(214, 7).Value = 130
$ws.Cells.Item(214, 9).Value = 53
$ws.Cells.Item(215, 7).Value = 10
$ws.Cells.Item(215, 9).Value = 82
$ws.Cells.Item(216, 7).Value = 12
$ws.Cells.Item(216, 9).Value = 80
$ws.Cells.Item(217, 7).Value = 86
$ws.Cells.Item(217, 9).Value = 6
$ws.Cells.Item(218, 7).Value = 49
$ws.Cells.Item(218, 9).Value = 42
$ws.Cells.Item(219, 7).Value = 54
$ws.Cells.Item(219, 9).Value = 37
$ws.Cells.Item(220, 7).Value = 40
$ws.Cells.Item(220, 9).Value = 51
$ws.Cells.Item(221, 7).Value = 36
$ws.Cells.Item(221, 9).Value = 56
$ws.Cells.Item(222, 7).Value = 9
$ws.Cells.Item(222, 9).Value = 83
$ws.Cells.Item(223, 7).Value = 88
$ws.Cells.Item(223, 9).Value = 4
$ws.Cells.Item(224, 7).Value = 44
$ws.Cells.Item(224, 9).Value = 47
$ws.Cells.Item(225, 7).Value = 100
$ws.Cells.Item(225, 9).Value = 84
$ws.Cells.Item(226, 7).Value = 100
$ws.Cells.Item(226, 9).Value = 84
$ws.Cells.Item(227, 7).Value = 142
$ws.Cells.Item(227, 9).Value = 41
$ws.Cells.Item(228, 7).Value = 155
$ws.Cells.Item(228, 9).Value = 28
$ws.Cells.Item(229, 7).Value = 137
$ws.Cells.Item(229, 9).Value = 46
$ws.Cells.Item(230, 7).Value = 35
$ws.Cells.Item(230, 9).Value = 148
$ws.Cells.Item(231, 7).Value = 22
$ws.Cells.Item(231, 9).Value = 161
$ws.Cells.Item(232, 9).Value = 58
$ws.Cells.Item(233, 7).Value = 157
$ws.Cells.Item(233, 9).Value = 26
$ws.Cells.Item(234, 7).Value = 37
$ws.Cells.Item(234, 9).Value = 146
$ws.Cells.Item(235, 7).Value = 32
$ws.Cells.Item(235, 9).Value = 151
$ws.Cells.Item(236, 7).Value = 32
$ws.Cells.Item(236, 9).Value = 151
$ws.Cells.Item(237, 7).Value = 156
$ws.Cells.Item(237, 9).Value = 27
$ws.Cells.Item(238, 7).Value = 7
$ws.Cells.Item(238, 9).Value = 176
$ws.Cells.Item(239, 7).Value = 32
$ws.Cells.Item(239, 9).Value = 60
$ws.Cells.Item(240, 7).Value = 32
$ws.Cells.Item(240, 9).Value = 60
$ws.Cells.Item(241, 7).Value = 37
$ws.Cells.Item(241, 9).Value = 146
$ws.Cells.Item(242, 7).Value = 166
$ws.Cells.Item(242, 9).Value = 18
$ws.Cells.Item(243, 7).Value = 68
$ws.Cells.Item(243, 9).Value = 24
$ws.Cells.Item(244, 7).Value = 26
$ws.Cells.Item(244, 9).Value = 157
$ws.Cells.Item(245, 7).Value = 159
$ws.Cells.Item(245, 9).Value = 24
$ws.Cells.Item(246, 7).Value = 136
$ws.Cells.Item(246, 9).Value = 47
$ws.Cells.Item(247, 9).Value = 45
$ws.Cells.Item(248, 7).Value = 9
$ws.Cells.Item(249, 7).Value = 135
$ws.Cells.Item(249, 9).Value = 47
$ws.Cells.Item(250, 7).Value = 10
$ws.Cells.Item(251, 7).Value = 78
$ws.Cells.Item(251, 9).Value = 14
$ws.Cells.Item(252, 7).Value = 36
$ws.Cells.Item(252, 9).Value = 56
$ws.Cells.Item(253, 9).Value = 81
$ws.Cells.Item(254, 7).Value = 14
$ws.Cells.Item(254, 9).Value = 169
$ws.Cells.Item(255, 7).Value = 15
$ws.Cells.Item(255, 9).Value = 77
$ws.Cells.Item(256, 7).Value = 79
$ws.Cells.Item(256, 9).Value = 105
$ws.Cells.Item(257, 7).Value = 173
$ws.Cells.Item(257, 9).Value = 11
$ws.Cells.Item(258, 9).Value = 18
$ws.Cells.Item(259, 9).Value = 250
$ws.Cells.Item(260, 7).Value = 163
$ws.Cells.Item(260, 9).Value = 203
$ws.Cells.Item(261, 9).Value = 87
$ws.Cells.Item(262, 9).Value = 67
